# Updated symbol list on Sun Dec 25 23:04:03 UTC 2022 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Hora"/hour (column G) values for the
# crypto snapshot rows (2-51) on the active sheet. Column G moves from hour 22
# to hour 23 for every row; column D is only touched for the coins whose price
# actually moved in this run's data pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row number -> @{ column letter = new value }
$updates = @{
    2 = @{ "D" = "242.84"; "G" = "23" }
    3 = @{ "G" = "23" }
    4 = @{ "D" = "5.394"; "G" = "23" }
    5 = @{ "D" = "0.05976"; "G" = "23" }
    6 = @{ "D" = "3.398"; "G" = "23" }
    7 = @{ "D" = "6.482"; "G" = "23" }
    8 = @{ "D" = "0.8127"; "G" = "23" }
    9 = @{ "D" = "0.9087"; "G" = "23" }
    10 = @{ "D" = "0.1412"; "G" = "23" }
    11 = @{ "D" = "0.07406"; "G" = "23" }
    12 = @{ "D" = "0.03335"; "G" = "23" }
    13 = @{ "D" = "0.03062"; "G" = "23" }
    14 = @{ "D" = "0.09337"; "G" = "23" }
    15 = @{ "D" = "3.850"; "G" = "23" }
    16 = @{ "D" = "0.001577"; "G" = "23" }
    17 = @{ "D" = "0.04638"; "G" = "23" }
    18 = @{ "D" = "0.0005939"; "G" = "23" }
    19 = @{ "D" = "0.006135"; "G" = "23" }
    20 = @{ "D" = "0.005029"; "G" = "23" }
    21 = @{ "D" = "0.0009821"; "G" = "23" }
    22 = @{ "D" = "0.00007798"; "G" = "23" }
    23 = @{ "D" = "0.0002900"; "G" = "23" }
    24 = @{ "D" = "3.624"; "G" = "23" }
    25 = @{ "D" = "2.163"; "G" = "23" }
    26 = @{ "G" = "23" }
    27 = @{ "D" = "0.1297"; "G" = "23" }
    28 = @{ "G" = "23" }
    29 = @{ "G" = "23" }
    30 = @{ "G" = "23" }
    31 = @{ "G" = "23" }
    32 = @{ "G" = "23" }
    33 = @{ "G" = "23" }
    34 = @{ "G" = "23" }
    35 = @{ "G" = "23" }
    36 = @{ "G" = "23" }
    37 = @{ "G" = "23" }
    38 = @{ "G" = "23" }
    39 = @{ "G" = "23" }
    40 = @{ "D" = "0.03879"; "G" = "23" }
    41 = @{ "D" = "0.006179"; "G" = "23" }
    42 = @{ "D" = "0.1070"; "G" = "23" }
    43 = @{ "D" = "0.002800"; "G" = "23" }
    44 = @{ "D" = "0.007214"; "G" = "23" }
    45 = @{ "D" = "0.00005193"; "G" = "23" }
    46 = @{ "G" = "23" }
    47 = @{ "D" = "0.0005799"; "G" = "23" }
    48 = @{ "G" = "23" }
    49 = @{ "G" = "23" }
    50 = @{ "D" = "0.00002100"; "G" = "23" }
    51 = @{ "D" = "0.0002000"; "G" = "23" }
}

foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $cell = $ws.Range("$col$row")
        # Leading apostrophe forces the new value to stay text (matches the
        # workbook's existing inline-string cells instead of becoming numeric).
        $cell.Value = "'" + $updates[$row][$col]
        # Reset to the default style so the quote-prefix tweak above does not
        # leave behind a new/extra cell style vs. the original formatting.
        $cell.Style = "Normal"
    }
}
